$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.738.17'
$ws.Range('E2').Value = '  -1.44%  '

$ws.Range('D3').Value = '3.495.53'
$ws.Range('E3').Value = '  -3.74%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.53'
$ws.Range('E5').Value = '  -4.40%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '193.66'
$ws.Range('E6').Value = '  -2.84%  '

$ws.Range('E7').Value = '  -2.20%  '

$ws.Range('D8').Value = '3.483.08'
$ws.Range('E8').Value = '  -3.72%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.04%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.206'
$ws.Range('E10').Value = '  -6.83%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.620'
$ws.Range('E11').Value = '  -4.48%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '51.73'
$ws.Range('E12').Value = '  -4.27%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000288'
$ws.Range('E13').Value = '  -6.12%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.16'
$ws.Range('E14').Value = '  -4.21%  '

$ws.Range('D15').Value = '4.027.54'
$ws.Range('E15').Value = '  -4.33%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '649.04'
$ws.Range('E16').Value = '  -4.59%  '

$ws.Range('D17').Value = '69.570.16'
$ws.Range('E17').Value = '  -1.89%  '

$ws.Range('D18').Value = '3.477.92'
$ws.Range('E18').Value = '  -4.95%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.35'
$ws.Range('E19').Value = '  -5.22%  '

$ws.Range('E20').Value = '  -1.77%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.26'
$ws.Range('E21').Value = '  -3.99%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.951'
$ws.Range('E22').Value = '  -4.99%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.11'
$ws.Range('E23').Value = '  -2.89%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.26'
$ws.Range('E24').Value = '  -2.24%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '99.34'
$ws.Range('E25').Value = '  -6.14%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.29'
$ws.Range('E26').Value = '  -7.29%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.89'
$ws.Range('E27').Value = '  -3.90%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.07'
$ws.Range('E28').Value = '  -3.52%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.37'
$ws.Range('E29').Value = '  -4.97%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.72'
$ws.Range('E30').Value = '  -4.39%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.28'
$ws.Range('E31').Value = '  -8.25%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.76'
$ws.Range('E32').Value = '  -5.77%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.66'
$ws.Range('E33').Value = '  -4.55%  '

$ws.Range('E34').Value = '  -4.85%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '61.10'
$ws.Range('E35').Value = '  -3.50%  '

$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '529.50'
$ws.Range('E36').Value = '  +4.93%  '

$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '3.728.10'
$ws.Range('E37').Value = '  -5.80%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.08%  '

$ws.Range('D39').Value = '0.0₃0790'
$ws.Range('E39').Value = '  -8.92%  '

$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.55'
$ws.Range('E40').Value = '  +0.19%  '

$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.94'
$ws.Range('E41').Value = '  -3.44%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.376'
$ws.Range('E42').Value = '  -3.18%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.56'
$ws.Range('E43').Value = '  +70.17%  '

$ws.Range('E44').Value = '  -2.79%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '34.41'
$ws.Range('E45').Value = '  -6.50%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0444'
$ws.Range('E46').Value = '  -3.46%  '

$ws.Range('E47').Value = '  -3.84%  '

$ws.Range('E48').Value = '  -9.02%  '

$ws.Range('E49').Value = '  -4.21%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.19'
$ws.Range('E51').Value = '  -5.76%  '
